$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicated "Predicted Eg_RF" / "Predicted Eg" columns (D and E):
# keep the header cells' formatting (style) but clear their content, and
# fully clear the data cells below them (rows 2:35).
$ws.Range("D1:E35").ClearContents()

# New custom width for column D (Random Forest predictions column)
$ws.Columns("D").ColumnWidth = 15.8

# Move the active selection to L10 (matches the saved selection state)
$ws.Range("L10").Select()
